$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FRONT_STACK_LIST")

# Insert a new row at position 5 (pushes jest/jotai/react/typescript/vue/webpack down by one)
$ws.Rows.Item(5).Insert()

# Fix up formatting of the newly inserted row so column A keeps the bordered/bold/centered style
# used by every other "stack" cell in this column (the plain Insert() does not copy the border).
$newCell = $ws.Range("A5")
$newCell.Font.Bold = $true
$newCell.HorizontalAlignment = -4108   # xlCenter
$newCell.VerticalAlignment = -4160     # xlTop
$newCell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$newCell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$newCell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$newCell.Borders.Item(10).LineStyle = 1  # xlEdgeRight

# Populate the new row: "javascript" with count 1
$ws.Range("A5").Value = "javascript"
$ws.Range("B5").Value = 1

# The "react" count increases from 3 to 4 as part of this data refresh
$ws.Range("B8").Value = 4
